# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets
# to reflect the refreshed data pulled from bilibili at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 25
$ws1.Range("F6").Value = 5271
$ws1.Range("F7").Value = 178
$ws1.Range("F8").Value = 91
$ws1.Range("F10").Value = 358

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 25
$ws4.Range("F9").Value = 5271
$ws4.Range("F10").Value = 178
$ws4.Range("F11").Value = 91
$ws4.Range("F14").Value = 358
